# Insert a new data row at row 64 (pushing existing rows 64-161 down to 65-162)
# and populate it with a new Caqui (Mankaki, Primera) price entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(64).Insert()

$ws.Cells.Item(64, 1).Value  = 6
$ws.Cells.Item(64, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(64, 3).Value  = "Metropolitana"
$ws.Cells.Item(64, 4).Value  = 45028
$ws.Cells.Item(64, 5).Value  = 13
$ws.Cells.Item(64, 6).Value  = "Fruta"
$ws.Cells.Item(64, 7).Value  = 100107
$ws.Cells.Item(64, 8).Value  = "Otros"
$ws.Cells.Item(64, 9).Value  = 100107001
$ws.Cells.Item(64, 10).Value = "Caqui"
$ws.Cells.Item(64, 11).Value = "Mankaki"
$ws.Cells.Item(64, 12).Value = "Primera"
$ws.Cells.Item(64, 13).Value = 100
$ws.Cells.Item(64, 14).Value = 16000
$ws.Cells.Item(64, 15).Value = 17000
$ws.Cells.Item(64, 16).Value = 16500
$ws.Cells.Item(64, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(64, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(64, 19).Value = 1650
$ws.Cells.Item(64, 20).Value = 10
